$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Comment on D1: expand the "1:" target description to also cover 2/3 ---
$c = $ws.Range("D1").Comment
$c.Text("Author:`n0：自己`n1:  乙方全体`n2:  敌人`n3:  敌方全体")

# --- Header J3: auras#rate_id -> auras#rate_id_lv ---
$ws.Range("J3").Value = "auras#rate_id_lv"

# --- J column data rows: normalize id patterns (exporter fix) ---
# (order matches the upstream commit's shared-string insertion order)
$ws.Range("J5").Value  = "20_2001_1|80_2002_1"
$ws.Range("J6").Value  = "20_2001_2|80_2002_2"

$ws.Range("J10").Value = "90_2003_1"
$ws.Range("J11").Value = "90_2003_2"
$ws.Range("J12").Value = "90_2003_3"
$ws.Range("J13").Value = "90_2003_4"
$ws.Range("J14").Value = "90_2003_5"

$ws.Range("J15").Value = "90_2003_1"
$ws.Range("J16").Value = "90_2003_2"
$ws.Range("J17").Value = "90_2003_3"
$ws.Range("J18").Value = "90_2003_4"
$ws.Range("J19").Value = "90_2003_5"

$ws.Range("J20").Value = "90_2003_1"
$ws.Range("J21").Value = "90_2003_2"
$ws.Range("J22").Value = "90_2003_3"
$ws.Range("J23").Value = "90_2003_4"
$ws.Range("J24").Value = "90_2003_5"

$ws.Range("J25").Value = "90_2003_1"
$ws.Range("J26").Value = "90_2003_2"
$ws.Range("J27").Value = "90_2003_3"
$ws.Range("J28").Value = "90_2003_4"
$ws.Range("J29").Value = "90_2003_5"

$ws.Range("J7").Value  = "20_2001_1|80_2002_4"
$ws.Range("J8").Value  = "20_2001_1|80_2002_4"
$ws.Range("J9").Value  = "20_2001_1|80_2002_5"

# --- Update the active selection left on the sheet ---
$null = $ws.Range("J7").Select()
